# Add worksheet "Ark9" as the last sheet (after "Ark8") and make it active,
# matching the new <sheet name="Ark9" .../> entry + activeTab change in workbook.xml.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Ark9"

# Header row (row 2), reusing the same shared strings as the other sheets.
$ws.Range("A2").Value = "initial"
$ws.Range("B2").Value = "distance"
$ws.Range("C2").Value = "final"
$ws.Range("D2").Value = "distance"

# Data rows 3-102: column A = initial distance, column B = final distance.
$data = @"
3	19.681000000000001	4.7068099999999999
4	7.2486899999999999	3.7665099999999998
5	12.0977	3492.47
6	14.4453	6.2738100000000001
7	18.603899999999999	5.0152000000000001
8	16.677099999999999	6.8936599999999997
9	18.3157	7.00631
10	12.977	4.1956499999999997
11	14.9846	41.105400000000003
12	17.848500000000001	5.2284100000000002
13	15.5824	12.452199999999999
14	15.7865	5.6434699999999998
15	13.893700000000001	22.733899999999998
16	19.9543	6.6444400000000003
17	18.941500000000001	6.9756299999999998
18	19.6629	33.0413
19	11.982100000000001	4.9547600000000003
20	19.337199999999999	10.2332
21	4.5030700000000001	1.4656100000000001
22	14.4718	56.599800000000002
23	10.851000000000001	255.70099999999999
24	13.364100000000001	23.469100000000001
25	16.158999999999999	19.7639
26	19.143699999999999	12.083399999999999
27	18.476500000000001	4.6040000000000001
28	17.917200000000001	49.410699999999999
29	12.2585	39.158299999999997
30	17.616099999999999	6.8544099999999997
31	16.656300000000002	11.2644
32	19.308599999999998	7.0667799999999996
33	17.760899999999999	6.52996
34	18.191600000000001	28.668600000000001
35	14.968999999999999	62.945099999999996
36	13.582800000000001	4.2919900000000002
37	17.247	5.3877100000000002
38	12.285500000000001	40.470599999999997
39	19.0854	2.8316699999999999
40	9.3107699999999998	107.08499999999999
41	19.0733	5.75793
42	15.1214	124.812
43	7.46957	5.6046100000000001
44	12.3087	3.5590999999999999
45	9.8798499999999994	6.3836599999999999
46	14.5817	6.3193799999999998
47	17.730399999999999	7.2483500000000003
48	12.443099999999999	122.523
49	15.176	4.9315100000000003
50	16.822099999999999	106.18
51	12.8066	5.1778199999999996
52	14.272	7.3676899999999996
53	13.081099999999999	2642.07
54	6.2914599999999998	1.0716300000000001
55	17.182200000000002	47.158999999999999
56	15.3965	4.4881200000000003
57	10.865500000000001	28.571000000000002
58	6.6582499999999998	38.154400000000003
59	14.055	6.4270199999999997
60	15.892899999999999	5.2398600000000002
61	19.744700000000002	7.9282700000000004
62	10.986800000000001	3.70994
63	18.695900000000002	5.5580299999999996
64	8.1631199999999993	3.8704100000000001
65	17.677499999999998	49.976599999999998
66	16.760100000000001	297.13099999999997
67	16.262799999999999	4.4987000000000004
68	16.4634	33.192500000000003
69	19.045500000000001	7.5705
70	10.575100000000001	83.69
71	17.932600000000001	209.999
72	16.802299999999999	4.6197999999999997
73	12.851100000000001	2.9603600000000001
74	17.397300000000001	6.0581100000000001
75	13.985799999999999	3.7327300000000001
76	14.4496	3.9341900000000001
77	13.061	6.9902800000000003
78	16.309699999999999	5.33439
79	8.00624	2.6654200000000001
80	14.492699999999999	166.06800000000001
81	16.044799999999999	6.1206699999999996
82	14.7384	5.0501699999999996
83	3.9839799999999999	4.26891
84	11.5158	220.285
85	19.808599999999998	6.05532
86	16.898900000000001	5.0494599999999998
87	18.830200000000001	5.26485
88	19.125299999999999	5.8116899999999996
89	14.7654	5.53878
90	15.951599999999999	143.041
91	9.6877099999999992	4.9323699999999997
92	17.671500000000002	209.267
93	11.144299999999999	5.14154
94	14.476800000000001	24.183499999999999
95	15.232699999999999	164.14400000000001
96	15.8706	50.3489
97	12.898300000000001	3.76993
98	17.137699999999999	259.84300000000002
99	13.122299999999999	3.32524
100	17.953099999999999	5.9745600000000003
101	17.761700000000001	3.4159299999999999
102	15.257899999999999	4.0402199999999997
"@

foreach ($line in ($data -split "`n")) {
    $parts = $line -split "`t"
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = [double]$parts[1]
    $ws.Cells.Item($r, 2).Value = [double]$parts[2]
}

# Match the author's final selection on the new sheet.
$ws.Range("O17").Select() | Out-Null
